$d = $word.ActiveDocument

function Replace-WithLineBreak($findText, $replaceText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
}

Replace-WithLineBreak "cells.2.Mechanism" "cells.^l2.Mechanism"
Replace-WithLineBreak "techniques.3.Bioreactors" "techniques.^l3.Bioreactors"
Replace-WithLineBreak "selection.4.Application" "selection.^l4.Application"
